# Update prism templates for capture run time
#
# Target sheet: "EC_Prism_Template".
#  - Rename header "Prism No Tag File Name" -> "Prism File Name"
#  - Insert a new "Runtime (ms)" column right after "Prob for EC with mu=0"
#    (pushes "EC Exists with mu" / "Prob for EC with mu" one column right)
#  - Append a trailing "Runtime (ms)" column at the end of the table
#  - Drop the spare blank row below the table
#  - Re-fit column widths to the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EC_Prism_Template")

# --- 1. Rename "Prism No Tag File Name" -> "Prism File Name" (col E header) ---
$ws.Cells.Item(2, 5).Value = "Prism File Name"

# --- 2. Insert a new column before the old column H ("EC Exists with mu") ---
#         so it becomes the new "Runtime (ms)" column; everything after slides
#         one column to the right.
$ws.Columns.Item(8).Insert()
$ws.Cells.Item(2, 8).Value = "Runtime (ms)"

# Merge its 2-row span first, then (re)apply the header look (centered,
# thin box border) -- doing the formatting after the merge keeps both
# merged cells on the same uniform thin-box style, just like the sheet's
# pre-existing header merges (e.g. A2:A3), instead of the split
# top/bottom-only border pair a pre-merge format would otherwise leave
# behind on the two physical cells.
$ws.Range("H2:H3").Merge()
$newHdr = $ws.Range("H2:H3")
$newHdr.Borders.LineStyle = 1
$newHdr.HorizontalAlignment = -4108   # xlCenter
$newHdr.VerticalAlignment = -4108     # xlCenter

# --- 3. Append a trailing "Runtime (ms)" column after "Prob for EC with mu" ---
# (which is now column J after the insert above)
$ws.Cells.Item(2, 11).Value = "Runtime (ms)"
$ws.Range("K2:K3").Merge()
$trailHdr = $ws.Range("K2:K3")
$trailHdr.Borders.LineStyle = 1
$trailHdr.HorizontalAlignment = -4108
$trailHdr.VerticalAlignment = -4108

# Column K also needs a title-bar cell (K1) styled like the rest of row 1,
# even though the A1:J1 merge band itself stops at column J.
$row1K = $ws.Range("K1")
$row1K.Borders.LineStyle = 1
$row1K.HorizontalAlignment = -4108
$row1K.VerticalAlignment = -4108

# --- 4. Drop the now-superfluous blank row 4 below the table ---
$ws.Rows.Item(4).Delete()

# --- 5. Re-fit column widths for the new layout ---
$ws.Columns.Item(4).ColumnWidth = 6.6      # D  - Set of Subsets   (~7.46)
$ws.Columns.Item(5).ColumnWidth = 10.95    # E  - Prism File Name  (~11.9)
$ws.Columns.Item(6).ColumnWidth = 12.75    # F  - EC Exists mu=0   (~13.67)
$ws.Columns.Item(7).ColumnWidth = 12.75    # G  - Prob EC mu=0     (~13.67)
$ws.Columns.Item(8).ColumnWidth = 14.1     # H  - Runtime (ms)     (~14.99)
$ws.Columns.Item(9).ColumnWidth = 14.1     # I  - EC Exists mu     (~14.99)
$ws.Columns.Item(10).ColumnWidth = 12.25   # J  - Prob EC mu       (~13.12)
$ws.Columns.Item(11).ColumnWidth = 12.25   # K  - Runtime (ms)     (~13.12)
